$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp (col A) and production (col B) values for rows 2-97,
# per-row as the quarterly forecast model for PCSunEnergy was retrained
# (data shifted from 2025-06-17 to 2025-06-20).
$rowData = @{
    2 = @(45828.01041666666, 445)
    3 = @(45828.02083333334, 453)
    4 = @(45828.03125, 495)
    5 = @(45828.04166666666, 486)
    6 = @(45828.05208333334, 487)
    7 = @(45828.0625, 528)
    8 = @(45828.07291666666, 559)
    9 = @(45828.08333333334, 627)
    10 = @(45828.09375, 811)
    11 = @(45828.10416666666, 1028)
    12 = @(45828.11458333334, 1137)
    13 = @(45828.125, 1219)
    14 = @(45828.13541666666, 1396)
    15 = @(45828.14583333334, 1444)
    16 = @(45828.15625, 1481)
    17 = @(45828.16666666666, 1520)
    18 = @(45828.17708333334, 1768)
    19 = @(45828.1875, 1863)
    20 = @(45828.19791666666, 1919)
    21 = @(45828.20833333334, 1841)
    22 = @(45828.21875, 1969)
    23 = @(45828.22916666666, 2014)
    24 = @(45828.23958333334, 1991)
    25 = @(45828.25, 1938)
    26 = @(45828.26041666666, 1950)
    27 = @(45828.27083333334, 1950)
    28 = @(45828.28125, 2002)
    29 = @(45828.29166666666, 2027)
    30 = @(45828.30208333334, 2013)
    31 = @(45828.3125, 2043)
    32 = @(45828.32291666666, 2063)
    33 = @(45828.33333333334, 2020)
    34 = @(45828.34375, 1910)
    35 = @(45828.35416666666, 1864)
    36 = @(45828.36458333334, 1867)
    37 = @(45828.375, 1871)
    38 = @(45828.38541666666, 1817)
    39 = @(45828.39583333334, 1834)
    40 = @(45828.40625, 1821)
    41 = @(45828.41666666666, 1827)
    42 = @(45828.42708333334, 1686)
    43 = @(45828.4375, 1644)
    44 = @(45828.44791666666, 1539)
    45 = @(45828.45833333334, 1465)
    46 = @(45828.46875, 1329)
    47 = @(45828.47916666666, 1262)
    48 = @(45828.48958333334, $null)
    49 = @(45828.5, $null)
    50 = @(45828.51041666666, $null)
    51 = @(45828.52083333334, $null)
    52 = @(45828.53125, $null)
    53 = @(45828.54166666666, $null)
    54 = @(45828.55208333334, $null)
    55 = @(45828.5625, $null)
    56 = @(45828.57291666666, $null)
    57 = @(45828.58333333334, $null)
    58 = @(45828.59375, $null)
    59 = @(45828.60416666666, $null)
    60 = @(45828.61458333334, $null)
    61 = @(45828.625, $null)
    62 = @(45828.63541666666, $null)
    63 = @(45828.64583333334, $null)
    64 = @(45828.65625, $null)
    65 = @(45828.66666666666, $null)
    66 = @(45828.67708333334, $null)
    67 = @(45828.6875, $null)
    68 = @(45828.69791666666, $null)
    69 = @(45828.70833333334, $null)
    70 = @(45828.71875, $null)
    71 = @(45828.72916666666, $null)
    72 = @(45828.73958333334, $null)
    73 = @(45828.75, $null)
    74 = @(45828.76041666666, $null)
    75 = @(45828.77083333334, $null)
    76 = @(45828.78125, $null)
    77 = @(45828.79166666666, $null)
    78 = @(45828.80208333334, $null)
    79 = @(45828.8125, $null)
    80 = @(45828.82291666666, $null)
    81 = @(45828.83333333334, $null)
    82 = @(45828.84375, $null)
    83 = @(45828.85416666666, $null)
    84 = @(45828.86458333334, $null)
    85 = @(45828.875, $null)
    86 = @(45828.88541666666, $null)
    87 = @(45828.89583333334, $null)
    88 = @(45828.90625, $null)
    89 = @(45828.91666666666, $null)
    90 = @(45828.92708333334, $null)
    91 = @(45828.9375, $null)
    92 = @(45828.94791666666, $null)
    93 = @(45828.95833333334, $null)
    94 = @(45828.96875, $null)
    95 = @(45828.97916666666, $null)
    96 = @(45828.98958333334, $null)
    97 = @(45829, $null)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($row, 2).Value = $vals[1]
    }
}
